$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-06-30"

$ws.Range("A7").Value = "June (through 06-30)"

$ws.Range("D7").Value = 74
$ws.Range("E7").Value = 58
$ws.Range("F7").Value = 47
$ws.Range("G7").Value = 114
$ws.Range("H7").Value = 129
$ws.Range("I7").Value = 143

$ws.Range("D8").Value = 390
$ws.Range("E8").Value = 353
$ws.Range("F8").Value = 251
$ws.Range("G8").Value = 472
$ws.Range("H8").Value = 760
$ws.Range("I8").Value = 806
